$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge split runs of text into single runs (whitespace-only churn)
# ---------------------------------------------------------------------

$d.Content.Find.Execute("AUTORISATIONS ", $true, $false, $false, $false, $false, $true, 1, $false, "AUTORISATIONS ", 2) | Out-Null
$d.Content.Find.Execute("Statut 2 : ", $true, $false, $false, $false, $false, $true, 1, $false, "Statut 2 : ", 2) | Out-Null
$d.Content.Find.Execute("Statut 3 :", $true, $false, $false, $false, $false, $true, 1, $false, "Statut 3 :", 2) | Out-Null
$d.Content.Find.Execute("Total dépensé :", $true, $false, $false, $false, $false, $true, 1, $false, "Total dépensé :", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Strip the now-needless right-alignment from the three empty shaded
#    "value" cells that sit beside the Statut 1 / Statut 2 / Statut 3
#    labels (the paragraphs collapse down to a bare <w:p/>).
# ---------------------------------------------------------------------

function Clear-StatutValueCell($labelText) {
    $rng = $d.Content
    $rng.Find.Execute($labelText) | Out-Null
    $labelCell = $rng.Cells(1)
    $row = $labelCell.Row
    $valueCell = $row.Cells($labelCell.ColumnIndex + 1)
    $valueCell.Range.ParagraphFormat.Alignment = 0
}

Clear-StatutValueCell("Statut 1 :")
Clear-StatutValueCell("Statut 2 : ")
Clear-StatutValueCell("Statut 3 :")

# ---------------------------------------------------------------------
# 3) Replace the trailing "tab paragraph" at the end of the document with
#    a new "Les statuts :" heading, a 4-row recap table, and a closing
#    centered paragraph.
# ---------------------------------------------------------------------

$last = $d.Paragraphs.Last
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Les statuts</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t> :</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="Grilledutableau"/><w:tblW w:w="0" w:type="auto"/><w:tblInd w:w="2689" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3114"/><w:gridCol w:w="996"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="3114" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:t>Désignation</w:t></w:r><w:r><w:t> </w:t></w:r><w:r><w:t>autorisation</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="996" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Modifier</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3114" w:type="dxa"/></w:tcPr><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="996" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3114" w:type="dxa"/></w:tcPr><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="996" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3114" w:type="dxa"/></w:tcPr><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="996" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$last.Range.InsertXML($xml)

Write-Host "edit complete"
